$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 2.2929166666666667
$ws.Range("B7").NumberFormat = "[h]:mm:ss"
$ws.Range("C7").Value = "Higehiro (Text with visuals, Japanese, New):36;  Harry Potter book 1 (Text-only, English, Familiar):32; iCarly (Audiovisual, English, Familiar):27;"
$ws.Range("D7").Value = "Watched children's shows from my childhood, and a bit of Minecraft let's plays, read Harry Potter and manga."

$ws.Range("B7").Select()
